# Applies the scheduled-runner price/profit refresh described in the commit.
# Each worksheet (crafting class) gets its currentAveragePrice* / Leve* columns
# (H, I, J, K, L, M, N) updated in place for the rows whose upstream market data changed.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 16505.63
$ws.Range("I98").Value = 17902.152
$ws.Range("K98").Value = 17902.152
$ws.Range("M98").Value = -16404.152
$ws.Range("H120").Value = 133388.25
$ws.Range("J120").Value = 133388.25
$ws.Range("L120").Value = 133388.25
$ws.Range("N120").Value = -143064.25
$ws.Range("H122").Value = 16505.63
$ws.Range("I122").Value = 17902.152
$ws.Range("K122").Value = 53706.45599999999
$ws.Range("M122").Value = -51256.45599999999
$ws.Range("H125").Value = 6237
$ws.Range("I125").Value = 16438.5
$ws.Range("K125").Value = 147946.5
$ws.Range("M125").Value = -145486.5
$ws.Range("H132").Value = 2742.7546
$ws.Range("I132").Value = 2603.1924
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 7809.5772
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -5279.5772
$ws.Range("N132").Value = -35060
$ws.Range("H137").Value = 9025.763000000001
$ws.Range("J137").Value = 2577.7058
$ws.Range("L137").Value = 7733.117400000001
$ws.Range("N137").Value = -12833.1174

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 19446.166
$ws.Range("I61").Value = 36071
$ws.Range("J61").Value = 7571.2856
$ws.Range("K61").Value = 36071
$ws.Range("L61").Value = 7571.2856
$ws.Range("M61").Value = -35859
$ws.Range("N61").Value = -7995.2856
$ws.Range("H74").Value = 7569.3335
$ws.Range("I74").Value = 34050
$ws.Range("K74").Value = 34050
$ws.Range("M74").Value = -33176
$ws.Range("H77").Value = 7569.3335
$ws.Range("I77").Value = 34050
$ws.Range("K77").Value = 170250
$ws.Range("M77").Value = -165882
$ws.Range("H132").Value = 3792.5789
$ws.Range("I132").Value = 2819.5417
$ws.Range("J132").Value = 5460.643
$ws.Range("K132").Value = 8458.625100000001
$ws.Range("L132").Value = 16381.929
$ws.Range("M132").Value = -5928.625100000001
$ws.Range("N132").Value = -21441.929
$ws.Range("H135").Value = 74853.28999999999
$ws.Range("J135").Value = 74853.28999999999
$ws.Range("L135").Value = 74853.28999999999
$ws.Range("N135").Value = -84993.28999999999
$ws.Range("H136").Value = 19446.166
$ws.Range("I136").Value = 36071
$ws.Range("J136").Value = 7571.2856
$ws.Range("K136").Value = 108213
$ws.Range("L136").Value = 22713.8568
$ws.Range("M136").Value = -105663
$ws.Range("N136").Value = -27813.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 49982.668
$ws.Range("J50").Value = 49982.668
$ws.Range("L50").Value = 49982.668
$ws.Range("N50").Value = -51130.668
$ws.Range("H52").Value = 39499
$ws.Range("J52").Value = 49248.5
$ws.Range("L52").Value = 49248.5
$ws.Range("N52").Value = -49774.5
$ws.Range("H121").Value = 39499
$ws.Range("J121").Value = 49248.5
$ws.Range("L121").Value = 49248.5
$ws.Range("N121").Value = -52742.5
$ws.Range("H134").Value = 13068
$ws.Range("I134").Value = 14723
$ws.Range("J134").Value = 6999.6665
$ws.Range("K134").Value = 44169
$ws.Range("L134").Value = 20998.9995
$ws.Range("M134").Value = -41634
$ws.Range("N134").Value = -26068.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 40735
$ws.Range("J70").Value = 40735
$ws.Range("L70").Value = 40735
$ws.Range("N70").Value = -41365
$ws.Range("H73").Value = 40735
$ws.Range("J73").Value = 40735
$ws.Range("L73").Value = 40735
$ws.Range("N73").Value = -42919
$ws.Range("H100").Value = 51750
$ws.Range("J100").Value = 61000
$ws.Range("L100").Value = 61000
$ws.Range("N100").Value = -63164
$ws.Range("H107").Value = 7207.4707
$ws.Range("I107").Value = 9312.846
$ws.Range("K107").Value = 9312.846
$ws.Range("M107").Value = -7392.846
$ws.Range("H119").Value = 40100
$ws.Range("J119").Value = 33500
$ws.Range("L119").Value = 33500
$ws.Range("N119").Value = -43176
$ws.Range("H134").Value = 4861.5
$ws.Range("I134").Value = 5790.12
$ws.Range("J134").Value = 1545
$ws.Range("K134").Value = 17370.36
$ws.Range("L134").Value = 4635
$ws.Range("M134").Value = -14835.36
$ws.Range("N134").Value = -9705

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3875.0557
$ws.Range("I103").Value = 5109.273
$ws.Range("J103").Value = 1935.5714
$ws.Range("K103").Value = 15327.819
$ws.Range("L103").Value = 5806.7142
$ws.Range("M103").Value = -14448.819
$ws.Range("N103").Value = -7564.7142
$ws.Range("H122").Value = 6056.0303
$ws.Range("I122").Value = 1726
$ws.Range("J122").Value = 6829.25
$ws.Range("K122").Value = 15534
$ws.Range("L122").Value = 61463.25
$ws.Range("M122").Value = -13084
$ws.Range("N122").Value = -66363.25
$ws.Range("H129").Value = 2372.3635
$ws.Range("J129").Value = 4203.2
$ws.Range("L129").Value = 12609.6
$ws.Range("N129").Value = -22609.6
$ws.Range("H140").Value = 1664.4445
$ws.Range("I140").Value = 1664.4445
$ws.Range("K140").Value = 4993.333500000001
$ws.Range("M140").Value = 186.6664999999994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 29950
$ws.Range("J45").Value = 29950
$ws.Range("L45").Value = 29950
$ws.Range("N45").Value = -31068
$ws.Range("H122").Value = 16014.429
$ws.Range("I122").Value = 14535.167
$ws.Range("J122").Value = 17123.875
$ws.Range("K122").Value = 43605.501
$ws.Range("L122").Value = 51371.625
$ws.Range("M122").Value = -41155.501
$ws.Range("N122").Value = -56271.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2264.1304
$ws.Range("J46").Value = 2526.3
$ws.Range("L46").Value = 2526.3
$ws.Range("N46").Value = -2902.3
$ws.Range("H122").Value = 5408.8057
$ws.Range("I122").Value = 4822.696
$ws.Range("J122").Value = 6445.769
$ws.Range("K122").Value = 14468.088
$ws.Range("L122").Value = 19337.307
$ws.Range("M122").Value = -12018.088
$ws.Range("N122").Value = -24237.307
$ws.Range("H130").Value = 20500
$ws.Range("I130").Value = 20500
$ws.Range("K130").Value = 20500
$ws.Range("M130").Value = -15480
$ws.Range("H132").Value = 623240.3
$ws.Range("I132").Value = 1065317.8
$ws.Range("J132").Value = 4332
$ws.Range("K132").Value = 3195953.4
$ws.Range("L132").Value = 12996
$ws.Range("M132").Value = -3193423.4
$ws.Range("N132").Value = -18056

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 69000
$ws.Range("J125").Value = 69000
$ws.Range("L125").Value = 69000
$ws.Range("N125").Value = -78840
$ws.Range("H132").Value = 6800.712
$ws.Range("I132").Value = 7485.8125
$ws.Range("K132").Value = 22457.4375
$ws.Range("M132").Value = -19927.4375
$ws.Range("H136").Value = 372679.28
$ws.Range("I136").Value = 554966.0600000001
$ws.Range("J136").Value = 8105.7144
$ws.Range("K136").Value = 1664898.18
$ws.Range("L136").Value = 24317.1432
$ws.Range("M136").Value = -1662348.18
$ws.Range("N136").Value = -29417.1432

